# Added Flow vs R1L to the cell data modeled by tissue slice code
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared string / header + summary block mirroring the existing
# B-column ("Kpl") summary block, but for the F-column ("Flow_Lac") data.
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = "=AVERAGE(F`$1:F`$3)"
$ws.Range("H39").Formula = "=AVERAGE(F`$4:F`$6)"
$ws.Range("I39").Formula = "=AVERAGE(F`$9:F`$11)"
$ws.Range("J39").Formula = "=AVERAGE(F`$13:F`$16)"

$ws.Range("G40").Formula = "=STDEV(F`$1:F`$3)/SQRT(COUNT(F`$1:F`$3))"
$ws.Range("H40").Formula = "=STDEV(F`$4:F`$6)/SQRT(COUNT(F`$4:F`$6))"
$ws.Range("I40").Formula = "=STDEV(F`$9:F`$11)/SQRT(COUNT(F`$9:F`$11))"
$ws.Range("J40").Formula = "=STDEV(F`$13:F`$16)/SQRT(COUNT(F`$13:F`$16))"

$ws.Columns.Item(6).Select()
